# Applies the edit described by the diff:
# - Inserts two new data rows (weekly price records) right before the former
#   row 485, shifting the existing rows 485-554 down to 487-556.
# - Dimension grows from A1:R554 to A1:R556.
# - Populates the two newly inserted rows with the new "Agricola del Norte
#   S.A. de Arica - Brocoli" records (Segunda / Tercera quality, 2023-05-31).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 485-486; this shifts old rows 485..554 down to
# 487..556 and inherits the number formatting (e.g. date format) of the row
# immediately above, matching the surrounding data rows.
$ws.Range("A485:A486").EntireRow.Insert()

# --- New row 485 (Segunda) ---
$ws.Cells.Item(485, 1).Value = 1
$ws.Cells.Item(485, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(485, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(485, 4).Value = 45077
$ws.Cells.Item(485, 5).Value = 15
$ws.Cells.Item(485, 6).Value = 100112023
$ws.Cells.Item(485, 7).Value = "Brócoli"
$ws.Cells.Item(485, 8).Value = "Sin especificar"
$ws.Cells.Item(485, 9).Value = "Segunda"
$ws.Cells.Item(485, 10).Value = 1200
$ws.Cells.Item(485, 11).Value = 700
$ws.Cells.Item(485, 12).Value = 800
$ws.Cells.Item(485, 13).Value = 750
$ws.Cells.Item(485, 14).Value = "`$/unidad"
$ws.Cells.Item(485, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(485, 16).Value = 750
$ws.Cells.Item(485, 17).Value = 1
$ws.Cells.Item(485, 18).Value = "Hortaliza"

# --- New row 486 (Tercera) ---
$ws.Cells.Item(486, 1).Value = 1
$ws.Cells.Item(486, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(486, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(486, 4).Value = 45077
$ws.Cells.Item(486, 5).Value = 15
$ws.Cells.Item(486, 6).Value = 100112023
$ws.Cells.Item(486, 7).Value = "Brócoli"
$ws.Cells.Item(486, 8).Value = "Sin especificar"
$ws.Cells.Item(486, 9).Value = "Tercera"
$ws.Cells.Item(486, 10).Value = 1400
$ws.Cells.Item(486, 11).Value = 500
$ws.Cells.Item(486, 12).Value = 600
$ws.Cells.Item(486, 13).Value = 550
$ws.Cells.Item(486, 14).Value = "`$/unidad"
$ws.Cells.Item(486, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(486, 16).Value = 550
$ws.Cells.Item(486, 17).Value = 1
$ws.Cells.Item(486, 18).Value = "Hortaliza"
